$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5069.5
$ws.Range("J40").Value = 5963
$ws.Range("L40").Value = 5963
$ws.Range("N40").Value = -6313

$ws.Range("H70").Value = 1700
$ws.Range("J70").Value = 1700
$ws.Range("L70").Value = 5100
$ws.Range("N70").Value = -5640

$ws.Range("H73").Value = 1700
$ws.Range("J73").Value = 1700
$ws.Range("L73").Value = 5100
$ws.Range("N73").Value = -6972

$ws.Range("H76").Value = 4458.4
$ws.Range("I76").Value = 4458.4
$ws.Range("K76").Value = 4458.4
$ws.Range("M76").Value = -4143.4

$ws.Range("H79").Value = 4458.4
$ws.Range("I79").Value = 4458.4
$ws.Range("K79").Value = 4458.4
$ws.Range("M79").Value = -3366.4

$ws.Range("H138").Value = 2667.25
$ws.Range("J138").Value = 3665
$ws.Range("L138").Value = 10995
$ws.Range("N138").Value = -21275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 8059.8184
$ws.Range("I102").Value = 3556.3333
$ws.Range("J102").Value = 9748.625
$ws.Range("K102").Value = 3556.3333
$ws.Range("L102").Value = 9748.625
$ws.Range("M102").Value = -1934.3333
$ws.Range("N102").Value = -12992.625

$ws.Range("H110").Value = 2201.5217
$ws.Range("I110").Value = 1199.1765
$ws.Range("J110").Value = 5041.5
$ws.Range("K110").Value = 1199.1765
$ws.Range("L110").Value = 5041.5
$ws.Range("M110").Value = 845.8235
$ws.Range("N110").Value = -9131.5

$ws.Range("H132").Value = 8395.799999999999

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = $null
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 42500
$ws.Range("J40").Value = 42500
$ws.Range("L40").Value = 42500
$ws.Range("N40").Value = -43030

$ws.Range("H99").Value = 2652.8333
$ws.Range("I99").Value = 3041.6
$ws.Range("K99").Value = 3041.6
$ws.Range("M99").Value = -1543.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -388

$ws.Range("H51").Value = 34548.625
$ws.Range("J51").Value = 41648.168
$ws.Range("L51").Value = 41648.168
$ws.Range("N51").Value = -43120.168

$ws.Range("H61").Value = 34548.625
$ws.Range("J61").Value = 41648.168
$ws.Range("L61").Value = 41648.168
$ws.Range("N61").Value = -42344.168

$ws.Range("H99").Value = 4986.3076
$ws.Range("I99").Value = 4577.7
$ws.Range("K99").Value = 4577.7
$ws.Range("M99").Value = -3079.7

$ws.Range("H105").Value = 1953.75
$ws.Range("I105").Value = 1506.3334
$ws.Range("K105").Value = 1506.3334
$ws.Range("M105").Value = 240.6666

$ws.Range("H126").Value = 4986.3076
$ws.Range("I126").Value = 4577.7
$ws.Range("K126").Value = 13733.1
$ws.Range("M126").Value = -11263.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1453642.5
$ws.Range("I4").Value = 2266800.5
$ws.Range("J4").Value = 515383.47
$ws.Range("K4").Value = 6800401.5
$ws.Range("L4").Value = 1546150.41
$ws.Range("M4").Value = -6800289.5
$ws.Range("N4").Value = -1546374.41

$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 850
$ws.Range("J68").Value = 1075
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 3225
$ws.Range("M68").Value = -1739
$ws.Range("N68").Value = -4847

$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 850
$ws.Range("J71").Value = 1075
$ws.Range("K71").Value = 7650
$ws.Range("L71").Value = 9675
$ws.Range("M71").Value = -3594
$ws.Range("N71").Value = -17787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2725
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 5150
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 5150
$ws.Range("M5").Value = -188
$ws.Range("N5").Value = -5374

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null

$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50518

$ws.Range("H80").Value = 4860.8
$ws.Range("I80").Value = 1432.6666
$ws.Range("K80").Value = 1432.6666
$ws.Range("M80").Value = -434.6666

$ws.Range("H83").Value = 4860.8
$ws.Range("I83").Value = 1432.6666
$ws.Range("K83").Value = 7163.333000000001
$ws.Range("M83").Value = -2171.333000000001

$ws.Range("H132").Value = 118470.9
$ws.Range("I132").Value = 165244.58
$ws.Range("K132").Value = 495733.74
$ws.Range("M132").Value = -493203.74

$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8404.666999999999
$ws.Range("I46").Value = 11530.75
$ws.Range("J46").Value = 7267.909
$ws.Range("K46").Value = 11530.75
$ws.Range("L46").Value = 7267.909
$ws.Range("M46").Value = -11342.75
$ws.Range("N46").Value = -7643.909

$ws.Range("H68").Value = 8916.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null

$ws.Range("H71").Value = 8916.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null

$ws.Range("H100").Value = 5776
$ws.Range("I100").Value = 2279.2727
$ws.Range("K100").Value = 2279.2727
$ws.Range("M100").Value = -1738.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2592.8333
$ws.Range("I2").Value = 2592.8333
$ws.Range("K2").Value = 2592.8333
$ws.Range("M2").Value = -2480.8333

$ws.Range("H54").Value = 53453.54
$ws.Range("J54").Value = 53453.54
$ws.Range("L54").Value = 53453.54
$ws.Range("N54").Value = -54493.54

$ws.Range("H81").Value = 2955.4443
$ws.Range("I81").Value = 2955.4443
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5910.8886
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -4849.8886

$ws.Range("H84").Value = 2955.4443
$ws.Range("I84").Value = 2955.4443
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 29554.443
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -24250.443

$ws.Range("H122").Value = 2995.5881
$ws.Range("I122").Value = 2104.5715
$ws.Range("J122").Value = 4434.923
$ws.Range("K122").Value = 6313.7145
$ws.Range("L122").Value = 13304.769
$ws.Range("M122").Value = -3863.7145
$ws.Range("N122").Value = -18204.769

$ws.Range("H126").Value = 2833.4546
$ws.Range("I126").Value = 2926.8
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 8780.400000000001
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -6310.400000000001
$ws.Range("N126").Value = -10640

$ws.Range("H132").Value = 3399.0222
$ws.Range("I132").Value = 3176.611
$ws.Range("K132").Value = 9529.832999999999
$ws.Range("M132").Value = -6999.832999999999

$ws.Range("H140").Value = 43258.168
$ws.Range("J140").Value = 43258.168
$ws.Range("L140").Value = 43258.168
$ws.Range("N140").Value = -53618.168

$ws.Range("H141").Value = 148352.89
$ws.Range("J141").Value = 148352.89
$ws.Range("L141").Value = 148352.89
$ws.Range("N141").Value = -158712.89
